# ============================================================
# Edit script: restructure PlayerPerformance workbook
#  1. Add 'Player Info' sheet before 'ODI Batting'
#  2. Update 'ODI Batting': rename MATCH_CARD_LINK -> MATCH_CODE,
#     replace URL values with bare match codes, clear stray empty
#     INNING_NUMBER cells on 'did not bat' rows
#  3. Add 'ODI Batting Extra' sheet after 'ODI Batting'
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- 1. Player Info sheet (new, placed before ODI Batting) ----
$wsInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$wsInfo.Name = "Player Info"

$wsInfo.Cells.Item(1, 1).Value = "ID"
$wsInfo.Cells.Item(1, 2).Value = "NAME"
$wsInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$wsInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"

$wsInfo.Cells.Item(2, 1).Value = "'3900"
$wsInfo.Cells.Item(2, 2).Value = "'Jonathan Marc Bairstow"
$wsInfo.Cells.Item(2, 3).Value = "'Right Handed"
$wsInfo.Cells.Item(2, 4).Value = "'Right Arm Medium"

$wsInfoHeader = $wsInfo.Range("A1:D1")
$wsInfoHeader.Font.Bold = $true
$wsInfoHeader.Borders.LineStyle = 1
$wsInfoHeader.HorizontalAlignment = -4108
$wsInfoHeader.VerticalAlignment = -4160

# ---- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ----
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

# Replace the full howstat URL with the bare numeric match code
$wsBatting.Cells.Item(2, 4).Value = "'3331"
$wsBatting.Cells.Item(3, 4).Value = "'3335"
$wsBatting.Cells.Item(4, 4).Value = "'3337"
$wsBatting.Cells.Item(5, 4).Value = "'3341"
$wsBatting.Cells.Item(6, 4).Value = "'3343"
$wsBatting.Cells.Item(7, 4).Value = "'3346"
$wsBatting.Cells.Item(8, 4).Value = "'3444"
$wsBatting.Cells.Item(9, 4).Value = "'3800"
$wsBatting.Cells.Item(10, 4).Value = "'3809"
$wsBatting.Cells.Item(11, 4).Value = "'3832"
$wsBatting.Cells.Item(12, 4).Value = "'3833"
$wsBatting.Cells.Item(13, 4).Value = "'3834"
$wsBatting.Cells.Item(14, 4).Value = "'3904"
$wsBatting.Cells.Item(15, 4).Value = "'3906"
$wsBatting.Cells.Item(16, 4).Value = "'3908"
$wsBatting.Cells.Item(17, 4).Value = "'3910"
$wsBatting.Cells.Item(18, 4).Value = "'3911"
$wsBatting.Cells.Item(19, 4).Value = "'3930"
$wsBatting.Cells.Item(20, 4).Value = "'3932"
$wsBatting.Cells.Item(21, 4).Value = "'3946"
$wsBatting.Cells.Item(22, 4).Value = "'3948"
$wsBatting.Cells.Item(23, 4).Value = "'3949"
$wsBatting.Cells.Item(24, 4).Value = "'3978"
$wsBatting.Cells.Item(25, 4).Value = "'4020"
$wsBatting.Cells.Item(26, 4).Value = "'4021"
$wsBatting.Cells.Item(27, 4).Value = "'4030"
$wsBatting.Cells.Item(28, 4).Value = "'4045"
$wsBatting.Cells.Item(29, 4).Value = "'4068"
$wsBatting.Cells.Item(30, 4).Value = "'4070"
$wsBatting.Cells.Item(31, 4).Value = "'4072"
$wsBatting.Cells.Item(32, 4).Value = "'4073"
$wsBatting.Cells.Item(33, 4).Value = "'4075"
$wsBatting.Cells.Item(34, 4).Value = "'4108"
$wsBatting.Cells.Item(35, 4).Value = "'4115"
$wsBatting.Cells.Item(36, 4).Value = "'4117"
$wsBatting.Cells.Item(37, 4).Value = "'4123"
$wsBatting.Cells.Item(38, 4).Value = "'4125"
$wsBatting.Cells.Item(39, 4).Value = "'4137"
$wsBatting.Cells.Item(40, 4).Value = "'4138"
$wsBatting.Cells.Item(41, 4).Value = "'4139"
$wsBatting.Cells.Item(42, 4).Value = "'4146"
$wsBatting.Cells.Item(43, 4).Value = "'4149"
$wsBatting.Cells.Item(44, 4).Value = "'4165"
$wsBatting.Cells.Item(45, 4).Value = "'4166"
$wsBatting.Cells.Item(46, 4).Value = "'4167"
$wsBatting.Cells.Item(47, 4).Value = "'4168"
$wsBatting.Cells.Item(48, 4).Value = "'4169"
$wsBatting.Cells.Item(49, 4).Value = "'4170"
$wsBatting.Cells.Item(50, 4).Value = "'4171"
$wsBatting.Cells.Item(51, 4).Value = "'4173"
$wsBatting.Cells.Item(52, 4).Value = "'4175"
$wsBatting.Cells.Item(53, 4).Value = "'4209"
$wsBatting.Cells.Item(54, 4).Value = "'4210"
$wsBatting.Cells.Item(55, 4).Value = "'4211"
$wsBatting.Cells.Item(56, 4).Value = "'4253"
$wsBatting.Cells.Item(57, 4).Value = "'4254"
$wsBatting.Cells.Item(58, 4).Value = "'4255"
$wsBatting.Cells.Item(59, 4).Value = "'4256"
$wsBatting.Cells.Item(60, 4).Value = "'4260"
$wsBatting.Cells.Item(61, 4).Value = "'4287"
$wsBatting.Cells.Item(62, 4).Value = "'4292"
$wsBatting.Cells.Item(63, 4).Value = "'4294"
$wsBatting.Cells.Item(64, 4).Value = "'4300"
$wsBatting.Cells.Item(65, 4).Value = "'4303"
$wsBatting.Cells.Item(66, 4).Value = "'4308"
$wsBatting.Cells.Item(67, 4).Value = "'4314"
$wsBatting.Cells.Item(68, 4).Value = "'4321"
$wsBatting.Cells.Item(69, 4).Value = "'4326"
$wsBatting.Cells.Item(70, 4).Value = "'4331"
$wsBatting.Cells.Item(71, 4).Value = "'4336"
$wsBatting.Cells.Item(72, 4).Value = "'4342"
$wsBatting.Cells.Item(73, 4).Value = "'4346"
$wsBatting.Cells.Item(74, 4).Value = "'4354"
$wsBatting.Cells.Item(75, 4).Value = "'4355"
$wsBatting.Cells.Item(76, 4).Value = "'4401"
$wsBatting.Cells.Item(77, 4).Value = "'4405"
$wsBatting.Cells.Item(78, 4).Value = "'4408"
$wsBatting.Cells.Item(79, 4).Value = "'4426"
$wsBatting.Cells.Item(80, 4).Value = "'4427"
$wsBatting.Cells.Item(81, 4).Value = "'4428"
$wsBatting.Cells.Item(82, 4).Value = "'4429"
$wsBatting.Cells.Item(83, 4).Value = "'4430"
$wsBatting.Cells.Item(84, 4).Value = "'4431"
$wsBatting.Cells.Item(85, 4).Value = "'4454"
$wsBatting.Cells.Item(86, 4).Value = "'4456"
$wsBatting.Cells.Item(87, 4).Value = "'4457"
$wsBatting.Cells.Item(88, 4).Value = "'4469"
$wsBatting.Cells.Item(89, 4).Value = "'4470"
$wsBatting.Cells.Item(90, 4).Value = "'4471"
$wsBatting.Cells.Item(91, 4).Value = "'4609"
$wsBatting.Cells.Item(92, 4).Value = "'4613"
$wsBatting.Cells.Item(93, 4).Value = "'4618"
$wsBatting.Cells.Item(94, 4).Value = "'4619"
$wsBatting.Cells.Item(95, 4).Value = "'4620"
$wsBatting.Cells.Item(96, 4).Value = "'4622"

# Clear the stray empty INNING_NUMBER placeholder cells (rows where the
# player did not bat) so the cell is truly blank, not an empty string
$wsBatting.Cells.Item(5, 2).ClearContents()
$wsBatting.Cells.Item(9, 2).ClearContents()
$wsBatting.Cells.Item(15, 2).ClearContents()
$wsBatting.Cells.Item(16, 2).ClearContents()
$wsBatting.Cells.Item(58, 2).ClearContents()
$wsBatting.Cells.Item(61, 2).ClearContents()
$wsBatting.Cells.Item(77, 2).ClearContents()
$wsBatting.Cells.Item(90, 2).ClearContents()
$wsBatting.Cells.Item(96, 2).ClearContents()

# ---- 3. ODI Batting Extra sheet (new, placed after ODI Batting) ----
$wsExtra = $wb.Worksheets.Add($null, $wb.Worksheets.Item("ODI Batting"))
$wsExtra.Name = "ODI Batting Extra"

$wsExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$wsExtra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$wsExtra.Cells.Item(1, 3).Value = "NUM_4"
$wsExtra.Cells.Item(1, 4).Value = "NUM_6"
$wsExtra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

$wsExtra.Cells.Item(2, 1).Value = "'4405"
$wsExtra.Cells.Item(2, 2).Value = ""
$wsExtra.Cells.Item(2, 3).Value = ""
$wsExtra.Cells.Item(2, 4).Value = ""
$wsExtra.Cells.Item(2, 5).Value = ""
$wsExtra.Cells.Item(2, 6).Value = "'NO"

$wsExtra.Cells.Item(3, 1).Value = "'4408"
$wsExtra.Cells.Item(3, 2).Value = ""
$wsExtra.Cells.Item(3, 3).Value = ""
$wsExtra.Cells.Item(3, 4).Value = ""
$wsExtra.Cells.Item(3, 5).Value = ""
$wsExtra.Cells.Item(3, 6).Value = "'NO"

$wsExtra.Cells.Item(4, 1).Value = "'4426"
$wsExtra.Cells.Item(4, 2).Value = ""
$wsExtra.Cells.Item(4, 3).Value = ""
$wsExtra.Cells.Item(4, 4).Value = ""
$wsExtra.Cells.Item(4, 5).Value = ""
$wsExtra.Cells.Item(4, 6).Value = "'NO"

$wsExtra.Cells.Item(5, 1).Value = "'4427"
$wsExtra.Cells.Item(5, 2).Value = 2
$wsExtra.Cells.Item(5, 3).Value = "'14"
$wsExtra.Cells.Item(5, 4).Value = "'2"
$wsExtra.Cells.Item(5, 5).Value = "'37.96%"
$wsExtra.Cells.Item(5, 6).Value = "'YES"

$wsExtra.Cells.Item(6, 1).Value = "'4428"
$wsExtra.Cells.Item(6, 2).Value = 2
$wsExtra.Cells.Item(6, 3).Value = "'0"
$wsExtra.Cells.Item(6, 4).Value = "'0"
$wsExtra.Cells.Item(6, 5).Value = "'1.22%"
$wsExtra.Cells.Item(6, 6).Value = "'NO"

$wsExtra.Cells.Item(7, 1).Value = "'4429"
$wsExtra.Cells.Item(7, 2).Value = 2
$wsExtra.Cells.Item(7, 3).Value = "'4"
$wsExtra.Cells.Item(7, 4).Value = "'4"
$wsExtra.Cells.Item(7, 5).Value = "'30.55%"
$wsExtra.Cells.Item(7, 6).Value = "'NO"

$wsExtra.Cells.Item(8, 1).Value = "'4430"
$wsExtra.Cells.Item(8, 2).Value = 2
$wsExtra.Cells.Item(8, 3).Value = "'0"
$wsExtra.Cells.Item(8, 4).Value = "'0"
$wsExtra.Cells.Item(8, 5).Value = ""
$wsExtra.Cells.Item(8, 6).Value = "'NO"

$wsExtra.Cells.Item(9, 1).Value = "'4431"
$wsExtra.Cells.Item(9, 2).Value = 2
$wsExtra.Cells.Item(9, 3).Value = "'12"
$wsExtra.Cells.Item(9, 4).Value = "'2"
$wsExtra.Cells.Item(9, 5).Value = "'37.09%"
$wsExtra.Cells.Item(9, 6).Value = "'NO"

$wsExtra.Cells.Item(10, 1).Value = "'4454"
$wsExtra.Cells.Item(10, 2).Value = 2
$wsExtra.Cells.Item(10, 3).Value = "'6"
$wsExtra.Cells.Item(10, 4).Value = "'7"
$wsExtra.Cells.Item(10, 5).Value = "'37.45%"
$wsExtra.Cells.Item(10, 6).Value = "'NO"

$wsExtra.Cells.Item(11, 1).Value = "'4456"
$wsExtra.Cells.Item(11, 2).Value = ""
$wsExtra.Cells.Item(11, 3).Value = ""
$wsExtra.Cells.Item(11, 4).Value = ""
$wsExtra.Cells.Item(11, 5).Value = ""
$wsExtra.Cells.Item(11, 6).Value = "'NO"

$wsExtra.Cells.Item(12, 1).Value = "'4457"
$wsExtra.Cells.Item(12, 2).Value = 2
$wsExtra.Cells.Item(12, 3).Value = "'0"
$wsExtra.Cells.Item(12, 4).Value = "'0"
$wsExtra.Cells.Item(12, 5).Value = "'0.31%"
$wsExtra.Cells.Item(12, 6).Value = "'NO"

$wsExtra.Cells.Item(13, 1).Value = "'4469"
$wsExtra.Cells.Item(13, 2).Value = 1
$wsExtra.Cells.Item(13, 3).Value = "'6"
$wsExtra.Cells.Item(13, 4).Value = "'1"
$wsExtra.Cells.Item(13, 5).Value = "'22.75%"
$wsExtra.Cells.Item(13, 6).Value = "'NO"

$wsExtra.Cells.Item(14, 1).Value = "'4470"
$wsExtra.Cells.Item(14, 2).Value = 2
$wsExtra.Cells.Item(14, 3).Value = "'3"
$wsExtra.Cells.Item(14, 4).Value = "'1"
$wsExtra.Cells.Item(14, 5).Value = "'11.89%"
$wsExtra.Cells.Item(14, 6).Value = "'NO"

$wsExtra.Cells.Item(15, 1).Value = "'4471"
$wsExtra.Cells.Item(15, 2).Value = ""
$wsExtra.Cells.Item(15, 3).Value = ""
$wsExtra.Cells.Item(15, 4).Value = ""
$wsExtra.Cells.Item(15, 5).Value = ""
$wsExtra.Cells.Item(15, 6).Value = "'NO"

$wsExtra.Cells.Item(16, 1).Value = "'4609"
$wsExtra.Cells.Item(16, 2).Value = 2
$wsExtra.Cells.Item(16, 3).Value = "'1"
$wsExtra.Cells.Item(16, 4).Value = "'0"
$wsExtra.Cells.Item(16, 5).Value = "'6.36%"
$wsExtra.Cells.Item(16, 6).Value = "'NO"

$wsExtra.Cells.Item(17, 1).Value = "'4613"
$wsExtra.Cells.Item(17, 2).Value = 2
$wsExtra.Cells.Item(17, 3).Value = "'6"
$wsExtra.Cells.Item(17, 4).Value = "'0"
$wsExtra.Cells.Item(17, 5).Value = "'15.45%"
$wsExtra.Cells.Item(17, 6).Value = "'NO"

$wsExtra.Cells.Item(18, 1).Value = "'4618"
$wsExtra.Cells.Item(18, 2).Value = 2
$wsExtra.Cells.Item(18, 3).Value = "'0"
$wsExtra.Cells.Item(18, 4).Value = "'0"
$wsExtra.Cells.Item(18, 5).Value = ""
$wsExtra.Cells.Item(18, 6).Value = "'NO"

$wsExtra.Cells.Item(19, 1).Value = "'4619"
$wsExtra.Cells.Item(19, 2).Value = ""
$wsExtra.Cells.Item(19, 3).Value = ""
$wsExtra.Cells.Item(19, 4).Value = ""
$wsExtra.Cells.Item(19, 5).Value = ""
$wsExtra.Cells.Item(19, 6).Value = "'NO"

$wsExtra.Cells.Item(20, 1).Value = "'4620"
$wsExtra.Cells.Item(20, 2).Value = 2
$wsExtra.Cells.Item(20, 3).Value = "'3"
$wsExtra.Cells.Item(20, 4).Value = "'0"
$wsExtra.Cells.Item(20, 5).Value = "'13.93%"
$wsExtra.Cells.Item(20, 6).Value = "'NO"

$wsExtra.Cells.Item(21, 1).Value = "'4622"
$wsExtra.Cells.Item(21, 2).Value = ""
$wsExtra.Cells.Item(21, 3).Value = ""
$wsExtra.Cells.Item(21, 4).Value = ""
$wsExtra.Cells.Item(21, 5).Value = ""
$wsExtra.Cells.Item(21, 6).Value = "'NO"

$wsExtraHeader = $wsExtra.Range("A1:F1")
$wsExtraHeader.Font.Bold = $true
$wsExtraHeader.Borders.LineStyle = 1
$wsExtraHeader.HorizontalAlignment = -4108
$wsExtraHeader.VerticalAlignment = -4160

# ---- Final sheet order check (Player Info, ODI Batting, ODI Batting Extra) ----
Write-Host "Sheets:"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host $i $wb.Worksheets.Item($i).Name
}
